$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6800
$ws.Range("J40").Value = 6333.3335
$ws.Range("L40").Value = 6333.3335
$ws.Range("N40").Value = -6683.3335
$ws.Range("H70").Value = 2555.2222
$ws.Range("I70").Value = 1857
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 5571
$ws.Range("L70").Value = 14997
$ws.Range("M70").Value = -5301
$ws.Range("N70").Value = -15537
$ws.Range("H73").Value = 2555.2222
$ws.Range("I73").Value = 1857
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 5571
$ws.Range("L73").Value = 14997
$ws.Range("M73").Value = -4635
$ws.Range("N73").Value = -16869

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 25998.5
$ws.Range("J10").Value = 50000
$ws.Range("L10").Value = 50000
$ws.Range("N10").Value = -50340
$ws.Range("H74").Value = 1106.4286
$ws.Range("I74").Value = 1020
$ws.Range("K74").Value = 1020
$ws.Range("M74").Value = -146
$ws.Range("H77").Value = 1106.4286
$ws.Range("I77").Value = 1020
$ws.Range("K77").Value = 5100
$ws.Range("M77").Value = -732
$ws.Range("H97").Value = 2402.25
$ws.Range("I97").Value = 2203
$ws.Range("K97").Value = 2203
$ws.Range("M97").Value = -1707

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1899.4
$ws.Range("I20").Value = 1899.4
$ws.Range("K20").Value = 1899.4
$ws.Range("M20").Value = -1652.4
$ws.Range("H99").Value = 2401.4285
$ws.Range("I99").Value = 2322
$ws.Range("J99").Value = 2600
$ws.Range("K99").Value = 2322
$ws.Range("L99").Value = 2600
$ws.Range("M99").Value = -824
$ws.Range("N99").Value = -5596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3929.2856
$ws.Range("I2").Value = 400
$ws.Range("K2").Value = 400
$ws.Range("M2").Value = -287
$ws.Range("H94").Value = 3550.3635
$ws.Range("J94").Value = 4305.75
$ws.Range("L94").Value = 4305.75
$ws.Range("N94").Value = -5207.75
$ws.Range("H107").Value = 986.4
$ws.Range("I107").Value = 1008
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1008
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 912
$ws.Range("N107").Value = -4740
$ws.Range("H122").Value = 1608.6666
$ws.Range("I122").Value = 1456
$ws.Range("K122").Value = 4368
$ws.Range("M122").Value = -1918
$ws.Range("H132").Value = 1300
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 666666800
$ws.Range("I4").Value = 666666800
$ws.Range("K4").Value = 2000000400
$ws.Range("M4").Value = -2000000288
$ws.Range("H8").Value = 875.8
$ws.Range("I8").Value = 875.8
$ws.Range("K8").Value = 2627.4
$ws.Range("M8").Value = -2488.4
$ws.Range("H34").Value = 566.2222
$ws.Range("I34").Value = 260
$ws.Range("J34").Value = 949
$ws.Range("K34").Value = 780
$ws.Range("L34").Value = 2847
$ws.Range("M34").Value = -696
$ws.Range("N34").Value = -3015
$ws.Range("H46").Value = 31.666666
$ws.Range("I46").Value = 31.666666
$ws.Range("K46").Value = 94.99999800000001
$ws.Range("M46").Value = -3.999998000000005
$ws.Range("H80").Value = 2979.8
$ws.Range("J80").Value = 2979.8
$ws.Range("L80").Value = 8939.400000000001
$ws.Range("N80").Value = -10811.4
$ws.Range("H83").Value = 2979.8
$ws.Range("J83").Value = 2979.8
$ws.Range("L83").Value = 26818.2
$ws.Range("N83").Value = -36178.2
$ws.Range("H112").Value = 43499.25
$ws.Range("J112").Value = 48570.57
$ws.Range("L112").Value = 145711.71
$ws.Range("N112").Value = -147927.71
$ws.Range("H129").Value = 1352
$ws.Range("I129").Value = 592.625
$ws.Range("J129").Value = 2364.5
$ws.Range("K129").Value = 1777.875
$ws.Range("L129").Value = 7093.5
$ws.Range("M129").Value = 3222.125
$ws.Range("N129").Value = -17093.5
$ws.Range("H131").Value = 986.53845
$ws.Range("J131").Value = 990
$ws.Range("L131").Value = 2970
$ws.Range("N131").Value = -13050

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13748
$ws.Range("I80").Value = 9997.333000000001
$ws.Range("K80").Value = 9997.333000000001
$ws.Range("M80").Value = -8999.333000000001
$ws.Range("H83").Value = 13748
$ws.Range("I83").Value = 9997.333000000001
$ws.Range("K83").Value = 49986.665
$ws.Range("M83").Value = -44994.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3752.1904
$ws.Range("I22").Value = 2870.4285
$ws.Range("J22").Value = 4193.0713
$ws.Range("K22").Value = 2870.4285
$ws.Range("L22").Value = 4193.0713
$ws.Range("M22").Value = -2575.4285
$ws.Range("N22").Value = -4783.0713
$ws.Range("H27").Value = 3752.1904
$ws.Range("I27").Value = 2870.4285
$ws.Range("J27").Value = 4193.0713
$ws.Range("K27").Value = 2870.4285
$ws.Range("L27").Value = 4193.0713
$ws.Range("M27").Value = -2763.4285
$ws.Range("N27").Value = -4407.0713
$ws.Range("H46").Value = 3500
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 3666.6667
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 3666.6667
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -4042.6667
$ws.Range("H55").Value = 635.2857
$ws.Range("I55").Value = 437
$ws.Range("K55").Value = 437
$ws.Range("M55").Value = -264

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H54").Value = 13000
$ws.Range("I54").Value = 2500
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 2500
$ws.Range("L54").Value = 20000
$ws.Range("M54").Value = -1980
$ws.Range("N54").Value = -21040
$ws.Range("H124").Value = 81999
$ws.Range("J124").Value = 81999
$ws.Range("L124").Value = 81999
$ws.Range("N124").Value = -91819
$ws.Range("H136").Value = 1120.6471
$ws.Range("I136").Value = 1118
$ws.Range("K136").Value = 3354
$ws.Range("M136").Value = -804
